$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50; existing rows 50-123 shift down to 51-124.
$ws.Rows("50:50").Insert()

# Populate the newly inserted row 50 with a new weekly price record.
# Columns A,B,C,E,F,G,H,I,K,L,M,N,O,P,Q,R mirror the (now shifted) row 51
# values, while D (Fecha) and J (Volumen) carry the new data point.
$ws.Range("A50").Value = 9
$ws.Range("B50").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C50").Value = "Metropolitana"
$ws.Range("D50").Value = 45070
$ws.Range("E50").Value = 13
$ws.Range("F50").Value = 100112005
$ws.Range("G50").Value = "Puerro"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 70
$ws.Range("K50").Value = 7000
$ws.Range("L50").Value = 7000
$ws.Range("M50").Value = 7000
$ws.Range("N50").Value = '$/paquete 20 unidades'
$ws.Range("O50").Value = "Provincia de Chacabuco"
$ws.Range("P50").Value = 350
$ws.Range("Q50").Value = 20
$ws.Range("R50").Value = "Hortaliza"
